$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 0
    3  = 2
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    9  = 4
    10 = 0
    11 = 1
    12 = 0
    13 = 0
    15 = 2
    16 = 4
    17 = 1
    18 = 2
    19 = 2
    21 = 0
    22 = 0
    23 = 1
    25 = 0
    26 = 3
    27 = 0
    28 = 1
    29 = 0
    30 = 1
    31 = 1
    32 = 1
    34 = 1
    35 = 1
    36 = 2
    37 = 1
    38 = 2
    39 = 1
    40 = 3
    41 = 1
    42 = 1
    43 = 2
    44 = 1
    45 = 2
    46 = 2
    47 = 4
    48 = 2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 7).Value = $values[$row]
}
